$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status text for both locale columns (E = zh-cn, F = de-de)
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E and F widen (to fit the longer status text) - matches the
# widened "Latest Target File" / "Latest Handback File" columns below.
$wsOverview.Range("E1").ColumnWidth = 29.2
$wsOverview.Range("F1").ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Latest Target File (I) now points at the handed-back markdown file, and
# becomes a hyperlink (same style as column A).
$wsZhCn.Range("I2").Value = "16a8ce33-3c70-4f0a-8592-fd990596d514.md"
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09cb99f396c51bdbb6296bc2f863e4cb78f65aa2/e2e/16a8ce33-3c70-4f0a-8592-fd990596d514.md", "", "", "16a8ce33-3c70-4f0a-8592-fd990596d514.md")

$wsZhCn.Range("I3").Value = "97d2c480-4af2-4b8f-984a-2d91f993666a.md"
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09cb99f396c51bdbb6296bc2f863e4cb78f65aa2/e2e/97d2c480-4af2-4b8f-984a-2d91f993666a.md", "", "", "97d2c480-4af2-4b8f-984a-2d91f993666a.md")

# Latest Handback File (J) now holds the generated handback xliff name.
$wsZhCn.Range("J2").Value = "16a8ce33-3c70-4f0a-8592-fd990596d514.d877e9dfdcc416726901bf18703191739e1f162f.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "97d2c480-4af2-4b8f-984a-2d91f993666a.90b6328456069e6c366ab9635a397214d2beb591.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-24 20:30:12"
$wsZhCn.Range("K3").Value = "2016-08-24 20:30:12"

# Latest Target File / Latest Handback File columns widen to fit full names.
$wsZhCn.Range("C1").ColumnWidth = 29.2
$wsZhCn.Range("I1").ColumnWidth = 39.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("I2").Value = "16a8ce33-3c70-4f0a-8592-fd990596d514.md"
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09cb99f396c51bdbb6296bc2f863e4cb78f65aa2/e2e/16a8ce33-3c70-4f0a-8592-fd990596d514.md", "", "", "16a8ce33-3c70-4f0a-8592-fd990596d514.md")

$wsDeDe.Range("I3").Value = "97d2c480-4af2-4b8f-984a-2d91f993666a.md"
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09cb99f396c51bdbb6296bc2f863e4cb78f65aa2/e2e/97d2c480-4af2-4b8f-984a-2d91f993666a.md", "", "", "97d2c480-4af2-4b8f-984a-2d91f993666a.md")

$wsDeDe.Range("J2").Value = "16a8ce33-3c70-4f0a-8592-fd990596d514.d877e9dfdcc416726901bf18703191739e1f162f.de-de.xlf"
$wsDeDe.Range("J3").Value = "97d2c480-4af2-4b8f-984a-2d91f993666a.90b6328456069e6c366ab9635a397214d2beb591.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-24 20:30:27"
$wsDeDe.Range("K3").Value = "2016-08-24 20:30:27"

$wsDeDe.Range("C1").ColumnWidth = 29.2
$wsDeDe.Range("I1").ColumnWidth = 39.17
$wsDeDe.Range("J1").ColumnWidth = 39.17

Write-Output "Generated handback report."
